# Update weekly price data rows (4-14) for Hortaliza, Femacal de La Calera - Pepino dulce.
# The sheet's data blocks (grouped by Fecha) are being rotated to reflect a newer weekly
# snapshot. Column I (Calidad) stays anchored to its row position; columns D, J, K, L, M,
# N, P, Q move together with their data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Fecha, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $Unidad, $PrecioKg, $KgUnidades) {
    $ws.Cells.Item($Row, 4).Value  = $Fecha            # D - Fecha
    $ws.Cells.Item($Row, 10).Value = $Volumen          # J - Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMinimo     # K - Precio minimo
    $ws.Cells.Item($Row, 12).Value = $PrecioMaximo     # L - Precio maximo
    $ws.Cells.Item($Row, 13).Value = $PrecioPromedio   # M - Precio promedio ponderado
    $ws.Cells.Item($Row, 14).Value = $Unidad           # N - Unidad de comercializacion
    $ws.Cells.Item($Row, 16).Value = $PrecioKg         # P - Precio $/Kg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades       # Q - Kg o Unidades
}

Set-Row 4  44424 75 18000 18000 18000 "`$/caja 15 kilos"    1200 15
Set-Row 5  44424 50 12000 12000 12000 "`$/caja 15 kilos"    800  15
Set-Row 6  44235 80 14000 14000 14000 "`$/bandeja 18 kilos" 778  18
Set-Row 7  44235 70 12000 12000 12000 "`$/bandeja 18 kilos" 667  18
Set-Row 8  44235 60 10000 10000 10000 "`$/bandeja 18 kilos" 556  18
Set-Row 9  44242 60 13000 13000 13000 "`$/bandeja 18 kilos" 722  18
Set-Row 10 44242 50 10000 10000 10000 "`$/bandeja 18 kilos" 556  18
Set-Row 11 44536 87 22000 22000 22000 "`$/bandeja 18 kilos" 1222 18
Set-Row 12 44536 80 20000 20000 20000 "`$/bandeja 18 kilos" 1111 18
Set-Row 13 44756 65 14000 14000 14000 "`$/caja 15 kilos"    933  15
Set-Row 14 44756 68 12000 12000 12000 "`$/caja 15 kilos"    800  15

# Calidad (column I) changes for rows 6-8 because the row-position labels
# (Primera/Segunda/Tercera) stay fixed while the underlying date block rotates.
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(7, 9).Value = "Segunda"
$ws.Cells.Item(8, 9).Value = "Tercera"
